# "Generate Report for Handback"
# Fills in the handback columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) on the per-locale sheets, flips the Status text
# from "Ready for handoff" to "Handed back: in sync with en-US", and widens
# the columns that now hold longer text/hyperlinks.

$wb = $excel.ActiveWorkbook

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c495eef0ee90cc11b7aae1bf2587b52b066b9f7/e2e/"

$md1 = "2e92161f-e707-4985-a61e-458d4157ca65.md"
$md2 = "fb00993f-6967-4fbc-85d1-e60b44d7aba1.md"

$zhXlf1 = "2e92161f-e707-4985-a61e-458d4157ca65.ad5a63ae3b75a15b7c1f2bbe0e4e7b5abfdde527.zh-cn.xlf"
$zhXlf2 = "fb00993f-6967-4fbc-85d1-e60b44d7aba1.e6c0f5bb19c953b6449ed4e52a9c489f7f45b18c.zh-cn.xlf"
$deXlf1 = "2e92161f-e707-4985-a61e-458d4157ca65.ad5a63ae3b75a15b7c1f2bbe0e4e7b5abfdde527.de-de.xlf"
$deXlf2 = "fb00993f-6967-4fbc-85d1-e60b44d7aba1.e6c0f5bb19c953b6449ed4e52a9c489f7f45b18c.de-de.xlf"

$handedBackStatus = "Handed back: in sync with en-US"

function Fill-LocaleSheet($sheetName, $xlf1, $xlf2, $handback1, $handback2) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status (column C) flips from "Ready for handoff" to the handed-back text.
    $ws.Range("C2").Value = $handedBackStatus
    $ws.Range("C3").Value = $handedBackStatus

    # Latest Target File (column I): the source markdown file, now hyperlinked.
    $ws.Range("I2").Value = $md1
    $ws.Range("I3").Value = $md2
    $ws.Hyperlinks.Add($ws.Range("I2"), ($mdUrlBase + $md1), "", "", $md1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), ($mdUrlBase + $md2), "", "", $md2) | Out-Null

    # Latest Handback File (column J): the locale-specific xlf.
    $ws.Range("J2").Value = $xlf1
    $ws.Range("J3").Value = $xlf2

    # Latest Handback DateTime (column K).
    $ws.Range("K2").Value = $handback1
    $ws.Range("K3").Value = $handback2

    # Columns I and J grew to fit the longer filenames/hyperlinks.
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
    # Column C (Status) grew to fit the longer status text.
    $ws.Columns.Item(3).ColumnWidth = 29.17
}

Fill-LocaleSheet "zh-cn" $zhXlf1 $zhXlf2 "2016-08-16 22:28:02" "2016-08-16 22:28:02"
Fill-LocaleSheet "de-de" $deXlf1 $deXlf2 "2016-08-16 22:28:12" "2016-08-16 22:28:12"

# Overview sheet: the zh-cn/de-de status columns (E/F) mirror the per-locale
# Status column, so they flip to the same "Handed back..." text and the
# columns widen to fit it.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $handedBackStatus
$overview.Range("F2").Value = $handedBackStatus
$overview.Range("E3").Value = $handedBackStatus
$overview.Range("F3").Value = $handedBackStatus
$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17
